$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns keep their textual nature (values like
# "64.291.27" or "10.00" must not be re-interpreted as numbers) while we
# overwrite them, then drop back to the default "Normal" style so we do not
# leave a stray number-format behind on cells we touch.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

# --- Row 2 (Bitcoin) ---
$ws.Range("D2").Value = "64.291.27"
$ws.Range("E2").Value = "  +0.61%  "

# --- Row 3 (Ethereum) ---
$ws.Range("D3").Value = "3.488.50"
$ws.Range("E3").Value = "  -0.32%  "

# --- Row 4 (TetherUSD) ---
$ws.Range("E4").Value = "  +0.01%  "

# --- Row 5 (BNB) ---
$ws.Range("D5").Value = "586.58"
$ws.Range("E5").Value = "  +0.37%  "

# --- Row 6 (Solana) ---
$ws.Range("D6").Value = "134.37"
$ws.Range("E6").Value = "  +1.86%  "

# --- Row 7 (LidoStakedEther) ---
$ws.Range("D7").Value = "3.488.54"
$ws.Range("E7").Value = "  -0.33%  "

# --- Row 8 (USDC) ---
$ws.Range("E8").Value = "  -0.04%  "

# --- Row 9 (XRP) ---
$ws.Range("E9").Value = "  -0.89%  "

# --- Row 10 (Dogecoin) ---
$ws.Range("E10").Value = "  -0.49%  "

# --- Row 11 (Toncoin) ---
$ws.Range("E11").Value = "  +1.54%  "

# --- Row 12 (Cardano) ---
$ws.Range("E12").Value = "  -2.81%  "

# --- Row 13 (WrappedliquidstakedEther2.0) ---
$ws.Range("D13").Value = "4.081.43"
$ws.Range("E13").Value = "  +0.13%  "

# --- Row 14 (TRON) ---
$ws.Range("E14").Value = "  +2.15%  "

# --- Row 15 (ShibaInu) ---
$ws.Range("E15").Value = "  +0.39%  "

# --- Row 16 (WrappedEther) ---
$ws.Range("D16").Value = "3.487.62"
$ws.Range("E16").Value = "  +0.00%  "

# --- Row 17 (WrappedBTC) ---
$ws.Range("D17").Value = "64.333.22"
$ws.Range("E17").Value = "  +0.60%  "

# --- Row 18 (Avalanche) ---
$ws.Range("D18").Value = "25.17"
$ws.Range("E18").Value = "  -9.60%  "

# --- Row 19 (Uniswap) ---
$ws.Range("D19").Value = "10.00"
$ws.Range("E19").Value = "  -0.35%  "

# --- Row 20 (Polkadot) ---
$ws.Range("E20").Value = "  +0.56%  "

# --- Row 21 (Chainlink) ---
$ws.Range("D21").Value = "13.69"

# --- Row 22 (BitcoinCash) ---
$ws.Range("D22").Value = "385.08"
$ws.Range("E22").Value = "  -1.68%  "

# --- Row 23 (Polygon) ---
$ws.Range("E23").Value = "  -2.30%  "

# --- Row 24 (WrappedeETH) ---
$ws.Range("D24").Value = "3.627.50"
$ws.Range("E24").Value = "  -0.15%  "

# --- Row 25 ---
$ws.Range("D25").Value = "74.16"
$ws.Range("E25").Value = "  +1.73%  "

# --- Row 26 ---
$ws.Range("E26").Value = "  +0.04%  "

# --- Row 27 ---
$ws.Range("D27").Value = "5.70"
$ws.Range("E27").Value = "  -0.65%  "

# --- Row 28 ---
$ws.Range("E28").Value = "  +0.73%  "

# --- Rows 29 & 30: Fetch.AI and Binance-PegBSC-USD swap places ---
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.03%  "

$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").Value = "1.54"
$ws.Range("E30").Value = "  -2.00%  "

# --- Row 31 ---
$ws.Range("E31").Value = "  -0.39%  "

# --- Row 32 ---
$ws.Range("E32").Value = "  -0.74%  "

# --- Row 33 ---
$ws.Range("D33").Value = "8.24"
$ws.Range("E33").Value = "  +0.30%  "

# --- Row 34 ---
$ws.Range("D34").Value = "3.509.14"
$ws.Range("E34").Value = "  +0.26%  "

# --- Row 36 ---
$ws.Range("E36").Value = "  +2.34%  "

# --- Row 37 ---
$ws.Range("D37").Value = "23.39"
$ws.Range("E37").Value = "  -1.83%  "

# --- Row 38 ---
$ws.Range("E38").Value = "  -1.76%  "

# --- Row 39 ---
$ws.Range("E39").Value = "  -1.81%  "

# --- Row 40 ---
$ws.Range("E40").Value = "  -2.48%  "

# --- Row 41 ---
$ws.Range("D41").Value = "161.86"
$ws.Range("E41").Value = "  -3.93%  "

# --- Row 42 ---
$ws.Range("D42").Value = "0.0777"
$ws.Range("E42").Value = "  -3.91%  "

# --- Row 43 ---
$ws.Range("E43").Value = "  -1.03%  "

# --- Row 44 ---
$ws.Range("E44").Value = "  +0.01%  "

# --- Row 45 ---
$ws.Range("D45").Value = "25.40"
$ws.Range("E45").Value = "  -3.39%  "

# --- Row 46 ---
$ws.Range("D46").Value = "41.74"
$ws.Range("E46").Value = "  -0.34%  "

# --- Row 47 ---
$ws.Range("E47").Value = "  +0.47%  "

# --- Row 48 ---
$ws.Range("D48").Value = "1.19"
$ws.Range("E48").Value = "  -0.48%  "

# --- Row 49 ---
$ws.Range("E49").Value = "  +0.26%  "

# --- Row 50 ---
$ws.Range("D50").Value = "2.467.66"
$ws.Range("E50").Value = "  +1.20%  "

# --- Row 51 ---
$ws.Range("D51").Value = "6.73"
$ws.Range("E51").Value = "  -2.26%  "

# Restore the default style on the touched columns so we don't leave a
# "Text" number format lingering on cells that didn't have one before.
$dataRange.Style = "Normal"
